$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "206.15") must stay as text,
# matching the original inline-string cell type; force text format first.
$forceTextCells = @("D5","D6","D8","D10","D15","D17","D18","D20","D23","D25","D26","D27","D30","D36","D37","D39","D40","D43","D44","D48","D49","D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.724.01"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "1.545.65"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "206.15"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "0.480"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "21.45"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "0.0581"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.767.26"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "1.541.49"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "26.721.48"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").Value = "61.22"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "212.77"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").Value = "  -4.73%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "152.38"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "6.48"
$ws.Range("E26").Value = "  -2.79%  "
$ws.Range("D27").Value = "14.87"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "0.0460"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "1.345.32"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "2.27"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").Value = "0.936"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +7.14%  "
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "62.66"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "1.680.61"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("D48").Value = "85.91"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "0.0506"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  +0.39%  "

# Restore default General number format now that the text values are committed.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "General"
}
